$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions data update).
# Some "Price" values are plain decimal numbers (e.g. "396.95"); Excel
# auto-converts numeric-looking text to a Number on assignment, so those
# cells are pre-formatted as Text to keep the literal string (matching the
# source data, which stores prices as text, e.g. "1.00", "0.580").

$ws.Range("D2").Value = "57.063.96"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "3.264.47"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.95"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.58"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  +4.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.41"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  +5.53%  "
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "3.769.08"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.30"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.01"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "3.252.72"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.91"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "56.919.30"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000109"
$ws.Range("E21").Value = "  +6.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.97"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "292.28"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.25"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.18"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.99"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.19"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.38"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.20"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.17"
$ws.Range("E34").Value = "  +10.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0486"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.30"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.19"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.93"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.87"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.283"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.70"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.24"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D49").Value = "2.149.65"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("E51").Value = "  -5.55%  "
